$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "247.85"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.557"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05630"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.478"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.073"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8020"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07308"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03196"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02992"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09264"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001670"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.971"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04682"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005915"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006268"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.001052"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.003832"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001501"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "UpBots"
$ws.Range("C22").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0004604"
$ws.Range("E22").Value = "21UpBotsUBXT"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.986"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "GateToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.398"
$ws.Range("E24").Value = "23GateTokenGT"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.113"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3311"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("B27").Value = "ProBitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1292"
$ws.Range("E27").Value = "26ProBitTokenPROBBestin24h"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007016"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1047"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002972"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008725"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005643"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02739"
